# Add a new column (H) holding the CO2/(CO+CO2) feed ratio to the Graaf data
# sheet. Inserting a blank column before the existing column H pushes the
# columns that were H..Q one slot to the right, becoming I..R.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank column before the current column H (8th column).
$ws.Columns.Item(8).Insert()

# Header for the newly inserted column.
$ws.Range("H1").Value = "CO2/(CO+CO2)"

# Row 2 gets its own formula...
$ws.Range("H2").Formula = "=F2/(E2+F2)"

# ...while rows 3-19 are written as one fill-down/shared formula group.
$ws.Range("H3:H19").Formula = "=F3/(E3+F3)"

# Give the whole new column the same plain/general look used elsewhere on
# the sheet (border only, no special number format) by copying A1's format.
$ws.Range("A1").Copy()
$ws.Range("H1:H19").PasteSpecial(-4122)

# Restore the selection the author left the sheet in.
$ws.Range("H2:H19").Select()
